# Updated cryptos list on Tue Feb 20 05:56:42 UTC 2024 with GitHub Actions
# Refreshes the Price (D) and Volume(1h) (E) columns on the active sheet
# with the latest scraped quotes. Numeric-looking price strings are
# written with a leading apostrophe to force Excel to keep them as text
# (matching the source data, which mixes thousand-separator dots with
# decimal dots, e.g. "51.934.64"); the style is then reset to "Normal"
# so no stray text-number-format is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.934.64"
$ws.Range("E2").Value = "  -0.27%  "

$ws.Range("D3").Value = "2.924.86"
$ws.Range("E3").Value = "  +0.96%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'358.11"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.09%  "

$ws.Range("D6").Value = "'110.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.93%  "

$ws.Range("D7").Value = "'0.569"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.98%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "'0.628"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.03%  "

$ws.Range("D10").Value = "'39.48"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.91%  "

$ws.Range("D11").Value = "'0.0885"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.15%  "

$ws.Range("E12").Value = "  +0.78%  "

$ws.Range("D13").Value = "'19.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.28%  "

$ws.Range("D14").Value = "'7.88"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.34%  "

$ws.Range("D15").Value = "3.390.32"
$ws.Range("E15").Value = "  +1.36%  "

$ws.Range("D16").Value = "2.918.84"
$ws.Range("E16").Value = "  +0.88%  "

$ws.Range("D17").Value = "'0.986"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.35%  "

$ws.Range("D18").Value = "51.940.57"
$ws.Range("E18").Value = "  -0.22%  "

$ws.Range("E19").Value = "  -0.24%  "

$ws.Range("D20").Value = "'7.58"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.47%  "

$ws.Range("D21").Value = "'14.08"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.54%  "

$ws.Range("D22").Value = "0.0₃0984"
$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("D23").Value = "'70.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.35%  "

$ws.Range("D24").Value = "'270.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.26%  "

$ws.Range("E25").Value = "  +1.86%  "

$ws.Range("E26").Value = "  +11.41%  "

$ws.Range("E27").Value = "  +2.32%  "

$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("E29").Value = "  +15.04%  "

$ws.Range("E30").Value = "  +14.22%  "

$ws.Range("D31").Value = "'10.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.94%  "

$ws.Range("D32").Value = "'38.77"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.13%  "

$ws.Range("E33").Value = "  -1.58%  "

$ws.Range("D34").Value = "'52.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.14%  "

$ws.Range("D35").Value = "'0.0444"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.88%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("E37").Value = "  -14.11%  "

$ws.Range("D38").Value = "'3.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.64%  "

$ws.Range("D39").Value = "'18.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("E40").Value = "  -1.19%  "

$ws.Range("D41").Value = "'2.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.97%  "

$ws.Range("D42").Value = "'0.121"
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").Value = "'23.28"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.90%  "

$ws.Range("D44").Value = "'119.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.20%  "

$ws.Range("D45").Value = "'2.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.64%  "

$ws.Range("E46").Value = "  -0.44%  "

$ws.Range("D47").Value = "'3.46"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.83%  "

$ws.Range("D48").Value = "2.139.48"
$ws.Range("E48").Value = "  -2.57%  "

$ws.Range("E49").Value = "  -8.33%  "

$ws.Range("E50").Value = "  +3.17%  "

$ws.Range("D51").Value = "'9.24"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.89%  "
